$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = -0.3455109999618997
$ws.Range("J3").Value = 0.2196478599460466
$ws.Range("K3").Value = -0.4328907103504786
$ws.Range("L3").Value = 2.801133796137951
